$d = $word.ActiveDocument

# --- Remove the stray "_GoBack" bookmark at the top of the document. ---
# Deleting it also causes Word to renumber the remaining w:id values of the
# IDX / IDX2 / IDX3 ... IDX28 bookmarks scattered through the rest of the
# document (each id shifts down by one).
$d.Bookmarks("_GoBack").Delete()

# --- Split the title into "5.2.1" and the old heading text. ---
# Original run text: "Stat 5100 Handout #14.c " + "– SAS: Logistic Regression with Polytomous Response"
# New text:          "5.2.1" + "– SAS: Logistic Regression with Polytomous Response"
$part1 = $d.Range(0, 24)
$part1.Text = "5.2.1"

# Touch the font of the remaining text so Word keeps it as its own run
# (identical formatting to the first run, just split into two <w:r> elements).
$part2 = $d.Range(5, 57)
$part2.Font.Name = "Times New Roman"
$part2.Font.NameBi = "Times New Roman"
